$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D values that look like plain decimals must be force-typed as Text
# (matching the source data, which stores prices as text strings) so Excel
# does not auto-convert them to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.481.61"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "2.605.72"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "572.81"
$ws.Range("E5").Value = "  +1.95%  "
$ws.Range("D6").Value = "142.70"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").Value = "2.631.23"
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("E10").Value = "  -2.79%  "
$ws.Range("D11").Value = "0.106"
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("E12").Value = "  -4.42%  "
$ws.Range("E13").Value = "  +2.30%  "
$ws.Range("D14").Value = "3.068.28"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").Value = "60.484.64"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("E17").Value = "  +2.56%  "
$ws.Range("D18").Value = "2.613.82"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "11.37"
$ws.Range("E19").Value = "  +9.01%  "
$ws.Range("D20").Value = "4.66"
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("D21").Value = "346.93"
$ws.Range("E21").Value = "  +2.41%  "
$ws.Range("D22").Value = "6.99"
$ws.Range("E22").Value = "  +6.98%  "
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").Value = "0.526"
$ws.Range("E24").Value = "  +11.43%  "
$ws.Range("D25").Value = "63.24"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("E28").Value = "  +3.48%  "
$ws.Range("D29").Value = "0.0₃0789"
$ws.Range("E29").Value = "  +0.97%  "
$ws.Range("E30").Value = "  +10.09%  "
$ws.Range("E31").Value = "  +3.08%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "161.61"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("D34").Value = "19.53"
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("D35").Value = "4.23"
$ws.Range("E35").Value = "  +3.94%  "
$ws.Range("D36").Value = "0.985"
$ws.Range("E36").Value = "  +10.07%  "
$ws.Range("E37").Value = "  +3.89%  "
$ws.Range("E38").Value = "  +7.55%  "
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("E40").Value = "  +4.24%  "
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("D42").Value = "294.89"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "137.46"
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").Value = "0.610"
$ws.Range("E45").Value = "  +2.20%  "
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("D47").Value = "19.76"
$ws.Range("E47").Value = "  +2.92%  "
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("E49").Value = "  +8.35%  "
$ws.Range("E50").Value = "  +1.63%  "
$ws.Range("D51").Value = "10.73"
$ws.Range("E51").Value = "  +0.88%  "

$ws.Range("D2:D51").Style = "Normal"
